$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.888.13'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.627.26'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.15'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.37'
$ws.Range('E8').Value = '  +8.36%  '
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('E10').Value = '  +1.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0915'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.860.05'
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.640.39'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.573'
$ws.Range('E14').Value = '  +6.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.90'
$ws.Range('E15').Value = '  +4.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.932.14'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.90'
$ws.Range('E17').Value = '  +17.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.69'
$ws.Range('E18').Value = '  +2.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.49'
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0706'
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.14'
$ws.Range('E22').Value = '  +3.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.63'
$ws.Range('E23').Value = '  +4.46%  '
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.53'
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.72'
$ws.Range('E26').Value = '  +2.53%  '
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.62'
$ws.Range('E28').Value = '  +3.28%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +3.17%  '
$ws.Range('E31').Value = '  +5.34%  '
$ws.Range('E32').Value = '  +4.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.22'
$ws.Range('E33').Value = '  +2.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.423.69'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('E35').Value = '  +6.84%  '
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('E40').Value = '  +3.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.831'
$ws.Range('E41').Value = '  +3.98%  '
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0497'
$ws.Range('E43').Value = '  +1.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '54.44'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.70'
$ws.Range('E45').Value = '  +5.81%  '
$ws.Range('E46').Value = '  +7.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.39'
$ws.Range('E48').Value = '  +2.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.767.69'
$ws.Range('E49').Value = '  +1.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '88.86'
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0108'
$ws.Range('E51').Value = '  +6.36%  '
